# Auto-generated Excel COM-interop script applying the Anima_Profits diff
# across all 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk 0: row 15 (date 44146)
$ws.Range("H15").Value = 2474.0337
$ws.Range("I15").Value = 2474.0337
$ws.Range("K15").Value = 7422.1011
$ws.Range("M15").Value = -7253.1011

# hunk 1: row 51 (date 5486)
$ws.Range("H51").Value = 1696.25
$ws.Range("I51").Value = 1949
$ws.Range("J51").Value = 1612
$ws.Range("K51").Value = 1949
$ws.Range("L51").Value = 1612
$ws.Range("M51").Value = -1465
$ws.Range("N51").Value = -2580

# hunk 2: row 94 (date 19905)
$ws.Range("H94").Value = 2900
$ws.Range("I94").Value = 2900
$ws.Range("K94").Value = 2900
$ws.Range("M94").Value = -2449

# hunk 3: row 98 (date 36237)
$ws.Range("H98").Value = 1138.7222
$ws.Range("I98").Value = 1060.9231
$ws.Range("J98").Value = 1341
$ws.Range("K98").Value = 1060.9231
$ws.Range("L98").Value = 1341
$ws.Range("M98").Value = 437.0769
$ws.Range("N98").Value = -4337

# hunk 4: row 122 (date 36237)
$ws.Range("H122").Value = 1138.7222
$ws.Range("I122").Value = 1060.9231
$ws.Range("J122").Value = 1341
$ws.Range("K122").Value = 3182.7693
$ws.Range("L122").Value = 4023
$ws.Range("M122").Value = -732.7692999999999
$ws.Range("N122").Value = -8923

# hunk 5: row 135 (date 44047)
$ws.Range("H135").Value = 1364.2593
$ws.Range("I135").Value = 1393.2
$ws.Range("K135").Value = 12538.8
$ws.Range("M135").Value = -10003.8

# hunk 6: row 137 (date 44013)
$ws.Range("H137").Value = 1227.76
$ws.Range("I137").Value = 1149.4375
$ws.Range("J137").Value = 1367
$ws.Range("K137").Value = 3448.3125
$ws.Range("L137").Value = 4101
$ws.Range("M137").Value = -898.3125
$ws.Range("N137").Value = -9201

# hunk 7: row 138 (date 44169)
$ws.Range("H138").Value = 2655.2144
$ws.Range("I138").Value = 3793.7
$ws.Range("J138").Value = 2407.7173
$ws.Range("K138").Value = 11381.1
$ws.Range("L138").Value = 7223.151899999999
$ws.Range("M138").Value = -6241.099999999999
$ws.Range("N138").Value = -17503.1519

$ws = $wb.Worksheets.Item("ARM")
# hunk 8: row 32 (date 44147)
$ws.Range("H32").Value = 570541.75
$ws.Range("I32").Value = 614941.7
$ws.Range("J32").Value = 144302.8
$ws.Range("K32").Value = 614941.7
$ws.Range("L32").Value = 144302.8
$ws.Range("M32").Value = -614654.7
$ws.Range("N32").Value = -144876.8

# hunk 9: row 45 (date 27714)
$ws.Range("H45").Value = 3266.7778
$ws.Range("I45").Value = 2406.889
$ws.Range("K45").Value = 2406.889
$ws.Range("M45").Value = -2029.889

# hunk 10: row 61 (date 43999)
$ws.Range("H61").Value = 10755079
$ws.Range("I61").Value = 25642496
$ws.Range("J61").Value = 3055.111
$ws.Range("K61").Value = 25642496
$ws.Range("L61").Value = 3055.111
$ws.Range("M61").Value = -25642284
$ws.Range("N61").Value = -3479.111

# hunk 11: row 74 (date 44000)
$ws.Range("H74").Value = 1146.8572
$ws.Range("I74").Value = 837.2857
$ws.Range("K74").Value = 837.2857
$ws.Range("M74").Value = 36.71429999999998

# hunk 12: row 77 (date 44000)
$ws.Range("H77").Value = 1146.8572
$ws.Range("I77").Value = 837.2857
$ws.Range("K77").Value = 4186.4285
$ws.Range("M77").Value = 181.5715

# hunk 13: row 123 (date 34107)
$ws.Range("H123").Value = 30122.818
$ws.Range("J123").Value = 30122.818
$ws.Range("L123").Value = 30122.818
$ws.Range("N123").Value = -39922.818

# hunk 14: row 132 (date 43997)
$ws.Range("H132").Value = 4882.7334
$ws.Range("I132").Value = 5737.2383
$ws.Range("J132").Value = 2888.889
$ws.Range("K132").Value = 17211.7149
$ws.Range("L132").Value = 8666.667000000001
$ws.Range("M132").Value = -14681.7149
$ws.Range("N132").Value = -13726.667

# hunk 15: row 136 (date 43999)
$ws.Range("H136").Value = 10755079
$ws.Range("I136").Value = 25642496
$ws.Range("J136").Value = 3055.111
$ws.Range("K136").Value = 76927488
$ws.Range("L136").Value = 9165.332999999999
$ws.Range("M136").Value = -76924938
$ws.Range("N136").Value = -14265.333

$ws = $wb.Worksheets.Item("BSM")
# hunk 16: row 109 (date 27096)
$ws.Range("H109").Value = 51333.332
$ws.Range("J109").Value = 51333.332
$ws.Range("L109").Value = 51333.332
$ws.Range("N109").Value = -54107.332

$ws = $wb.Worksheets.Item("CRP")
# hunk 17: row 31 (date 44023)
$ws.Range("H31").Value = 3045
$ws.Range("I31").Value = 1325.7368
$ws.Range("J31").Value = 3770.9111
$ws.Range("K31").Value = 1325.7368
$ws.Range("L31").Value = 3770.9111
$ws.Range("M31").Value = -1030.7368
$ws.Range("N31").Value = -4360.911099999999

# hunk 18: row 34 (date 44023)
$ws.Range("H34").Value = 3045
$ws.Range("I34").Value = 1325.7368
$ws.Range("J34").Value = 3770.9111
$ws.Range("K34").Value = 1325.7368
$ws.Range("L34").Value = 3770.9111
$ws.Range("M34").Value = -1123.7368
$ws.Range("N34").Value = -4174.911099999999

# hunk 19: row 58 (date 44021)
$ws.Range("H58").Value = 2571.6667
$ws.Range("I58").Value = 2368.4614
$ws.Range("J58").Value = 3100
$ws.Range("K58").Value = 2368.4614
$ws.Range("L58").Value = 3100
$ws.Range("M58").Value = -2165.4614
$ws.Range("N58").Value = -3506

# hunk 20: row 134 (date 44020)
$ws.Range("H134").Value = 963.03845
$ws.Range("I134").Value = 718.6818
$ws.Range("J134").Value = 2307
$ws.Range("K134").Value = 2156.0454
$ws.Range("L134").Value = 6921
$ws.Range("M134").Value = 378.9546
$ws.Range("N134").Value = -11991

# hunk 21: row 136 (date 44021)
$ws.Range("H136").Value = 2571.6667
$ws.Range("I136").Value = 2368.4614
$ws.Range("J136").Value = 3100
$ws.Range("K136").Value = 7105.3842
$ws.Range("L136").Value = 9300
$ws.Range("M136").Value = -4555.3842
$ws.Range("N136").Value = -14400

$ws = $wb.Worksheets.Item("CUL")
# hunk 22: row 68 (date 12895)
$ws.Range("H68").Value = 1351.7172
$ws.Range("I68").Value = 749.6786
$ws.Range("J68").Value = 1589.1409
$ws.Range("K68").Value = 2249.0358
$ws.Range("L68").Value = 4767.4227
$ws.Range("M68").Value = -1438.0358
$ws.Range("N68").Value = -6389.4227

# hunk 23: row 71 (date 12895)
$ws.Range("H71").Value = 1351.7172
$ws.Range("I71").Value = 749.6786
$ws.Range("J71").Value = 1589.1409
$ws.Range("K71").Value = 6747.1074
$ws.Range("L71").Value = 14302.2681
$ws.Range("M71").Value = -2691.1074
$ws.Range("N71").Value = -22414.2681

# hunk 24: row 105 (date 19814)
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# hunk 25: row 107 (date 27838)
$ws.Range("H107").Value = 1514.2297
$ws.Range("I107").Value = 317.26666
$ws.Range("K107").Value = 951.79998
$ws.Range("M107").Value = 968.20002

# hunk 26: row 122 (date 36078)
$ws.Range("H122").Value = 5786.7896
$ws.Range("I122").Value = 388.69232
$ws.Range("J122").Value = 17482.666
$ws.Range("K122").Value = 3498.23088
$ws.Range("L122").Value = 157343.994
$ws.Range("M122").Value = -1048.23088
$ws.Range("N122").Value = -162243.994

# hunk 27: row 129 (date 36054)
$ws.Range("H129").Value = 1194.7931
$ws.Range("I129").Value = 561.7273
$ws.Range("J129").Value = 1581.6666
$ws.Range("K129").Value = 1685.1819
$ws.Range("L129").Value = 4744.9998
$ws.Range("M129").Value = 3314.8181
$ws.Range("N129").Value = -14744.9998

# hunk 28: row 131 (date 36060)
$ws.Range("H131").Value = 1054.1305
$ws.Range("J131").Value = 1173.8235
$ws.Range("L131").Value = 3521.4705
$ws.Range("N131").Value = -13601.4705

# hunk 29: row 132 (date 43972)
$ws.Range("H132").Value = 3135.3662
$ws.Range("I132").Value = 2125.8484
$ws.Range("J132").Value = 4012.0527
$ws.Range("K132").Value = 19132.6356
$ws.Range("L132").Value = 36108.4743
$ws.Range("M132").Value = -16602.6356
$ws.Range("N132").Value = -41168.4743

# hunk 30: row 137 (date 44088)
$ws.Range("H137").Value = 7932.522
$ws.Range("I137").Value = 9703.429
$ws.Range("J137").Value = 5177.778
$ws.Range("K137").Value = 29110.287
$ws.Range("L137").Value = 15533.334
$ws.Range("M137").Value = -24010.287
$ws.Range("N137").Value = -25733.334

$ws = $wb.Worksheets.Item("GSM")
# hunk 31: row 122 (date 36182)
$ws.Range("H122").Value = 5907.909
$ws.Range("I122").Value = 7503.5
$ws.Range("J122").Value = 5553.3335
$ws.Range("K122").Value = 22510.5
$ws.Range("L122").Value = 16660.0005
$ws.Range("M122").Value = -20060.5
$ws.Range("N122").Value = -21560.0005

# hunk 32: row 126 (date 36184)
$ws.Range("H126").Value = 2254
$ws.Range("I126").Value = 2334
$ws.Range("J126").Value = 2014
$ws.Range("K126").Value = 7002
$ws.Range("L126").Value = 6042
$ws.Range("M126").Value = -4532
$ws.Range("N126").Value = -10982

# hunk 33: row 132 (date 44008)
$ws.Range("H132").Value = 2724.125
$ws.Range("I132").Value = 2622.6667
$ws.Range("J132").Value = 2785
$ws.Range("K132").Value = 7868.000100000001
$ws.Range("L132").Value = 8355
$ws.Range("M132").Value = -5338.000100000001
$ws.Range("N132").Value = -13415

$ws = $wb.Worksheets.Item("LTW")
# hunk 34: row 40 (date 36248)
$ws.Range("H40").Value = 333336670
$ws.Range("I40").Value = 1000000000
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 1000000000
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -999999864
$ws.Range("N40").Value = -5272

# hunk 35: row 69 (date 10671)
$ws.Range("H69").Value = 100163
$ws.Range("J69").Value = 100163
$ws.Range("L69").Value = 100163
$ws.Range("N69").Value = -101785

# hunk 36: row 72 (date 10671)
$ws.Range("H72").Value = 100163
$ws.Range("J72").Value = 100163
$ws.Range("L72").Value = 300489
$ws.Range("N72").Value = -308601

# hunk 37: row 94 (date 18067)
$ws.Range("H94").Value = 42886.668
$ws.Range("J94").Value = 42886.668
$ws.Range("L94").Value = 42886.668
$ws.Range("N94").Value = -44238.668

# hunk 38: row 109 (date 27209)
$ws.Range("H109").Value = 21370
$ws.Range("J109").Value = 21370
$ws.Range("L109").Value = 21370
$ws.Range("N109").Value = -24144

# hunk 39: row 132 (date 44058)
$ws.Range("H132").Value = 2535.8823
$ws.Range("I132").Value = 2107.1428
$ws.Range("K132").Value = 6321.428400000001
$ws.Range("M132").Value = -3791.428400000001

$ws = $wb.Worksheets.Item("WVR")
# hunk 40: row 75 (date 11957)
$ws.Range("H75").Value = 98130
$ws.Range("J75").Value = 98130
$ws.Range("L75").Value = 98130
$ws.Range("N75").Value = -100002

# hunk 41: row 76 (date 10896)
$ws.Range("H76").Value = 97782
$ws.Range("J76").Value = 97782
$ws.Range("L76").Value = 97782
$ws.Range("N76").Value = -98412

# hunk 42: row 78 (date 11957)
$ws.Range("H78").Value = 98130
$ws.Range("J78").Value = 98130
$ws.Range("L78").Value = 294390
$ws.Range("N78").Value = -303750

# hunk 43: row 79 (date 10896)
$ws.Range("H79").Value = 97782
$ws.Range("J79").Value = 97782
$ws.Range("L79").Value = 97782
$ws.Range("N79").Value = -99966

# hunk 44: row 122 (date 36208)
$ws.Range("H122").Value = 2740
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -13300

# hunk 45: row 123 (date 34127)
$ws.Range("H123").Value = 25121.857
$ws.Range("J123").Value = 25121.857
$ws.Range("L123").Value = 25121.857
$ws.Range("N123").Value = -34921.857

# hunk 46: row 126 (date 36210)
$ws.Range("H126").Value = 1267
$ws.Range("I126").Value = 1058.8
$ws.Range("J126").Value = 1527.25
$ws.Range("K126").Value = 3176.4
$ws.Range("L126").Value = 4581.75
$ws.Range("M126").Value = -706.3999999999996
$ws.Range("N126").Value = -9521.75

# hunk 47: row 136 (date 44031)
$ws.Range("H136").Value = 2311.224
$ws.Range("I136").Value = 2182.5813
$ws.Range("J136").Value = 2680
$ws.Range("K136").Value = 6547.743899999999
$ws.Range("L136").Value = 8040
$ws.Range("M136").Value = -3997.743899999999
$ws.Range("N136").Value = -13140
